$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.440.79'
$ws.Range("E2").Value = '  +3.10%  '
$ws.Range("D3").Value = '2.332.28'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '545.02'
$ws.Range("E5").Value = '  +1.35%  '
$ws.Range("D6").Value = '131.38'
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  -1.25%  '
$ws.Range("D9").Value = '2.329.59'
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("D14").Value = '23.65'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '2.747.40'
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = '60.415.79'
$ws.Range("E16").Value = '  +3.20%  '
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = '2.344.10'
$ws.Range("E18").Value = '  +1.74%  '
$ws.Range("D19").Value = '10.59'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").Value = '314.81'
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").Value = '64.10'
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '7.84'
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("E28").Value = '  +5.78%  '
$ws.Range("D29").Value = '1.20'
$ws.Range("E29").Value = '  +8.77%  '
$ws.Range("D30").Value = '173.01'
$ws.Range("E30").Value = '  +1.04%  '
$ws.Range("D31").Value = '1.73'
$ws.Range("E31").Value = '  +1.26%  '
$ws.Range("D32").Value = '0.0₃0733'
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("E33").Value = '  +1.68%  '
$ws.Range("E34").Value = '  +9.86%  '
$ws.Range("E35").Value = '  -1.24%  '
$ws.Range("D37").Value = '17.85'
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = '4.07'
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("D40").Value = '322.51'
$ws.Range("E40").Value = '  +10.91%  '
$ws.Range("D41").Value = '1.53'
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("D42").Value = '37.88'
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("D43").Value = '137.21'
$ws.Range("E43").Value = '  -2.82%  '
$ws.Range("D44").Value = '3.49'
$ws.Range("E44").Value = '  +0.84%  '
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("D46").Value = '19.14'
$ws.Range("E46").Value = '  +4.36%  '
$ws.Range("E47").Value = '  +1.27%  '
$ws.Range("D48").Value = '0.0495'
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("D49").Value = '0.0213'
$ws.Range("E49").Value = '  +1.29%  '
$ws.Range("D50").Value = '0.0₆0214'
$ws.Range("E50").Value = '  +15.12%  '
$ws.Range("D51").Value = '11.02'
